# Updated symbol list on Thu Dec 22 17:06:51 UTC 2022 with GitHub Actions
# Refresh the crypto pricing table: prices (col D), hour stamp (col G),
# and for several rows the coin/link/volume-label (cols B, C, E) which
# rotated to reflect the latest ranking snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row new values, keyed by column letter. Row numbers match the sheet.
$updates = @(
    @{ Row = 2; D = '242.12'; G = '17' },
    @{ Row = 3; D = '21.82'; G = '17' },
    @{ Row = 4; D = '5.374'; G = '17' },
    @{ Row = 5; D = '0.05697'; G = '17' },
    @{ Row = 6; D = '3.412'; G = '17' },
    @{ Row = 7; D = '6.300'; G = '17' },
    @{ Row = 8; D = '0.8056'; G = '17' },
    @{ Row = 9; D = '0.8436'; G = '17' },
    @{ Row = 10; B = 'WazirX'; C = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'; D = '0.1424'; E = '9WazirXWRX'; G = '17' },
    @{ Row = 11; B = 'MandalaExchangeToken'; C = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'; D = '0.07263'; E = '10MandalaExchangeTokenMDX'; G = '17' },
    @{ Row = 12; B = 'LiechtensteinCryptoassetsExchange'; C = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'; D = '0.03021'; E = '11LiechtensteinCryptoassetsExchangeLCX'; G = '17' },
    @{ Row = 13; B = 'BitrueCoin'; C = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'; D = '0.03147'; E = '12BitrueCoinBTR'; G = '17' },
    @{ Row = 14; B = 'BitMartToken'; C = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'; D = '0.09354'; E = '13BitMartTokenBMX'; G = '17' },
    @{ Row = 15; B = 'MCDex'; C = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'; D = '3.911'; E = '14MCDexMCB'; G = '17' },
    @{ Row = 16; B = 'BitForexToken'; C = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'; D = '0.001584'; E = '15BitForexTokenBF'; G = '17' },
    @{ Row = 17; B = 'CoinExToken'; C = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'; D = '0.04817'; E = '16CoinExTokenCET'; G = '17' },
    @{ Row = 18; B = 'One'; C = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'; D = '0.0005821'; E = '17OneONE'; G = '17' },
    @{ Row = 19; D = '0.006398'; G = '17' },
    @{ Row = 20; G = '17' },
    @{ Row = 21; D = '0.0009931'; G = '17' },
    @{ Row = 22; D = '0.0001500'; G = '17' },
    @{ Row = 23; D = '3.715'; G = '17' },
    @{ Row = 24; D = '2.148'; G = '17' },
    @{ Row = 25; G = '17' },
    @{ Row = 26; G = '17' },
    @{ Row = 27; D = '0.0004000'; G = '17' },
    @{ Row = 28; G = '17' },
    @{ Row = 29; G = '17' },
    @{ Row = 30; G = '17' },
    @{ Row = 31; G = '17' },
    @{ Row = 32; G = '17' },
    @{ Row = 33; G = '17' },
    @{ Row = 34; G = '17' },
    @{ Row = 35; G = '17' },
    @{ Row = 36; G = '17' },
    @{ Row = 37; G = '17' },
    @{ Row = 38; G = '17' },
    @{ Row = 39; G = '17' },
    @{ Row = 40; D = '0.03805'; G = '17' },
    @{ Row = 41; D = '0.006697'; G = '17' },
    @{ Row = 42; G = '17' },
    @{ Row = 43; D = '0.002621'; G = '17' },
    @{ Row = 44; D = '0.006865'; G = '17' },
    @{ Row = 45; D = '0.00005615'; G = '17' },
    @{ Row = 46; G = '17' },
    @{ Row = 47; D = '0.5801'; E = '46CoinbaseStockTokenCOINBestin24h'; G = '17' },
    @{ Row = 48; G = '17' },
    @{ Row = 49; D = '0.00002100'; G = '17' },
    @{ Row = 50; G = '17' },
    @{ Row = 51; G = '17' }

)

# Text columns (plain replace - these never look like numbers).
$textCols = @{ B = 2; C = 3; E = 5 }
# Numeric-looking columns that must still be stored as TEXT (General/auto
# type inference on assignment would otherwise turn "242.12" into a float
# and silently drop significant trailing zeros like "0.0001500").
$textForceCols = @{ D = 4; G = 7 }

foreach ($u in $updates) {
    $row = $u.Row

    foreach ($col in $textCols.Keys) {
        if ($u.ContainsKey($col)) {
            $ws.Cells.Item($row, $textCols[$col]).Value = $u[$col]
        }
    }

    foreach ($col in $textForceCols.Keys) {
        if ($u.ContainsKey($col)) {
            $cell = $ws.Cells.Item($row, $textForceCols[$col])
            $origStyle = $cell.Style
            $cell.NumberFormat = "@"
            $cell.Value = $u[$col]
            $cell.Style = $origStyle
        }
    }
}
